# A new weekly price record was inserted as row 360 ("Femacal de La Calera",
# Berenjena, Coquimbo) with fecha 45015. This pushes the previous rows
# 360-427 down to 361-428 (dimension grows from A1:R427 to A1:R428).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 360, shifting existing rows 360-427
# down to 361-428.
$ws.Rows(360).Insert()

# Populate the newly inserted row 360 with the new record's data.
$ws.Range("A360").Value = 3
$ws.Range("B360").Value = "Femacal de La Calera"
$ws.Range("C360").Value = "Coquimbo"
$ws.Range("D360").Value = 45015
$ws.Range("E360").Value = 5
$ws.Range("F360").Value = 100112001
$ws.Range("G360").Value = "Berenjena"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 55
$ws.Range("K360").Value = 8000
$ws.Range("L360").Value = 8000
$ws.Range("M360").Value = 8000
$ws.Range("N360").Value = "$/caja 60 unidades"
$ws.Range("O360").Value = "Región de Arica y Parinacota"
$ws.Range("P360").Value = 133
$ws.Range("Q360").Value = 60
$ws.Range("R360").Value = "Hortaliza"
